$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    -0.9716635481116656, 0.2805049688864443, -0.6760650065702858, -0.04573875985788579, -0.4726963440578849, 0.1769623127656392, 0.2080098490516883, 0.2409726375546795, 0.03734130806860806, 0.01272818083837047,
    0.8, 0.3, 0, 0.5, 0.1, 0.4, 0.5, 0.3, 0.1, 0.2,
    0.9, 0.4, 0, 0.5, 0.2132975555746283, 0.4704848369192122, 0.3964520361608751, 0.4230165803625844, 0.3531481090437554, 0.4376510431580233,
    0.6485846904589434, 0.5426591427781329, 0.2044328674106259, 0.4167513714050496, 0.5527780896650785, 0.4865322308196076, 0.5234718933384791, 0.535576229457463, 0.5969715903930113, 0.6598534028039785,
    0.6207578930310335, 0.2938136814264701, 0, 0.2964852479966615, 0.3, -0.1761138288871404, -0.02363626827304436, 0.6810025317521822, -6.9, 4.729401638091318,
    -0.8905127363963885, -2.04269378128221, 2.093024636165651, 1.098535546956398, -0.4717175472572421, 0.9401304606753627, 0.6392725048137464, 0.06357296580725347, -1.012166871044968, -0.5703626997413522,
    0.2394101325822788, -0.04072131480353391, -0.07958838003274593, 0.02912383308249389, -0.1311265493919933, -0.03907468377752821, 0.1771324545010202, 0.4946531409412387, 0.202428137729683, 0.208573386070384
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item(4 + $i, 2).Value = $newValues[$i]
}

$ws.Range("A74:B82").EntireRow.Delete()

Write-Output "Applied naive forecaster bugfix: updated B4:B73 and removed rows 74-82"
